try {
    $wb = $excel.ActiveWorkbook
    $ws = $wb.ActiveSheet
    $ws.Unprotect()

    # Update the confidentiality footer text (date 2021-05-19 -> 2021-05-20)
    $ws.Range("A80").Value2 = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-20 for illustrative purposes only and are subject to change."

    # Update Weight (D) and Percent Change (E) values for the holdings table (rows 2-77)
    $ws.Range("D2").Value2 = 0.06224687718455242
    $ws.Range("E2").Value2 = 0.0210121100328815
    $ws.Range("D3").Value2 = 0.03841325681691291
    $ws.Range("E3").Value2 = 0.004913670400395809
    $ws.Range("D4").Value2 = 0.031787035389135
    $ws.Range("E4").Value2 = 0.01382033563672258
    $ws.Range("D5").Value2 = 0.02936430683631882
    $ws.Range("E5").Value2 = 0.008237232289950436
    $ws.Range("D6").Value2 = 0.02699910664633258
    $ws.Range("E6").Value2 = 0.01560642747083407
    $ws.Range("D7").Value2 = 0.02553276731552215
    $ws.Range("E7").Value2 = -0.001737943020296617
    $ws.Range("D8").Value2 = 0.1938709789378046
    $ws.Range("E8").Value2 = 0.009870740305522929
    $ws.Range("D9").Value2 = 0.02493276075443614
    $ws.Range("E9").Value2 = 0.005820790216368721
    $ws.Range("D10").Value2 = 0.02273083493306028
    $ws.Range("E10").Value2 = 0.01024890190336758
    $ws.Range("D11").Value2 = 0.0220086400310875
    $ws.Range("E11").Value2 = 0.01224140041620769
    $ws.Range("D12").Value2 = 0.01877817681366088
    $ws.Range("E12").Value2 = 0.01234713770898543
    $ws.Range("D13").Value2 = 0.02028683039481731
    $ws.Range("E13").Value2 = -0.002382654276864571
    $ws.Range("D14").Value2 = 0.01725458646061264
    $ws.Range("E14").Value2 = 0.007242233657327901
    $ws.Range("D15").Value2 = 0.01626008271722781
    $ws.Range("E15").Value2 = 0.02359649122807017
    $ws.Range("D16").Value2 = 0.01480713497501003
    $ws.Range("E16").Value2 = 0.02483089305591224
    $ws.Range("D17").Value2 = 0.01461949425437457
    $ws.Range("E17").Value2 = 0.003523111612175889
    $ws.Range("D18").Value2 = 0.01440487225613085
    $ws.Range("E18").Value2 = 0.0123770886337069
    $ws.Range("D19").Value2 = 0.01366690856048981
    $ws.Range("E19").Value2 = 0.01600816352562262
    $ws.Range("D20").Value2 = 0.01308603850026252
    $ws.Range("E20").Value2 = -0.002373685995252561
    $ws.Range("D21").Value2 = 0.01135923673333353
    $ws.Range("E21").Value2 = 0.02348066298342544
    $ws.Range("D22").Value2 = 0.01331096172050059
    $ws.Range("E22").Value2 = 0.01087613293051382
    $ws.Range("D23").Value2 = 0.0114416465092883
    $ws.Range("E23").Value2 = 0.01861598980552936
    $ws.Range("D24").Value2 = 0.01297953971287483
    $ws.Range("E24").Value2 = -0.002100122100122115
    $ws.Range("D25").Value2 = 0.01121268204042168
    $ws.Range("E25").Value2 = 0.004611225950071685
    $ws.Range("D26").Value2 = 0.008867133412393891
    $ws.Range("E26").Value2 = 0.01501313649443259
    $ws.Range("D27").Value2 = 0.009543685976914804
    $ws.Range("E27").Value2 = 0.04071736964463635
    $ws.Range("D28").Value2 = 0.00992491043087476
    $ws.Range("E28").Value2 = 0.0118082889557769
    $ws.Range("D29").Value2 = 0.01004587054913913
    $ws.Range("E29").Value2 = 0.02031511900771021
    $ws.Range("D30").Value2 = 0.009676809461149417
    $ws.Range("E30").Value2 = 0.008516213560432195
    $ws.Range("D31").Value2 = 0.008554134744104176
    $ws.Range("E31").Value2 = 0.01065751445086716
    $ws.Range("D32").Value2 = 0.009885369586450312
    $ws.Range("E32").Value2 = 0.01474926253687325
    $ws.Range("D33").Value2 = 0.009246218381785793
    $ws.Range("E33").Value2 = -0.002283907238229244
    $ws.Range("D34").Value2 = 0.009014123926236278
    $ws.Range("E34").Value2 = 0.008860993169651055
    $ws.Range("D35").Value2 = 0.009200972245179859
    $ws.Range("E35").Value2 = -0.0008181544158806231
    $ws.Range("D36").Value2 = 0.008299773801013039
    $ws.Range("E36").Value2 = 0.01220140917683454
    $ws.Range("D37").Value2 = 0.008484244914881001
    $ws.Range("E37").Value2 = 0.00609414401793229
    $ws.Range("D38").Value2 = 0.006697299859538877
    $ws.Range("E38").Value2 = 0.04138714371916352
    $ws.Range("D39").Value2 = 0.008837180628441103
    $ws.Range("E39").Value2 = 0.007280944012051149
    $ws.Range("D40").Value2 = 0.007768706187119954
    $ws.Range("E40").Value2 = 0.00348837209302344
    $ws.Range("D41").Value2 = 0.007000393506680185
    $ws.Range("E41").Value2 = 0.03006429412297384
    $ws.Range("D42").Value2 = 0.007247464354206092
    $ws.Range("E42").Value2 = 0.01515383437930495
    $ws.Range("D43").Value2 = 0.008111182198347335
    $ws.Range("E43").Value2 = -0.0002540005080009511
    $ws.Range("D44").Value2 = 0.007375041026597597
    $ws.Range("E44").Value2 = 0.002331528279181816
    $ws.Range("D45").Value2 = 0.007241759062024609
    $ws.Range("E45").Value2 = 0.01825144983039717
    $ws.Range("D46").Value2 = 0.00792845436820149
    $ws.Range("E46").Value2 = -0.003837850803549925
    $ws.Range("D47").Value2 = 0.007546199792042099
    $ws.Range("E47").Value2 = 0.003780241935483986
    $ws.Range("D48").Value2 = 0.007169808988402568
    $ws.Range("E48").Value2 = 0.003558718861209842
    $ws.Range("D49").Value2 = 0.006517702016075929
    $ws.Range("E49").Value2 = 0.01586577915564868
    $ws.Range("D50").Value2 = 0.007775917042515995
    $ws.Range("E50").Value2 = -0.00156933079251198
    $ws.Range("D51").Value2 = 0.006569366606386029
    $ws.Range("E51").Value2 = -0.007936842994047399
    $ws.Range("D52").Value2 = 0.006644446666690967
    $ws.Range("E52").Value2 = 0.007584792253017225
    $ws.Range("D53").Value2 = 0.005484212109450959
    $ws.Range("E53").Value2 = -0.01582141309059393
    $ws.Range("D54").Value2 = 0.006228594258796173
    $ws.Range("E54").Value2 = 0.01536817464760065
    $ws.Range("D55").Value2 = 0.005426564886367221
    $ws.Range("E55").Value2 = 0.002190340597963081
    $ws.Range("D56").Value2 = 0.005672605611693694
    $ws.Range("E56").Value2 = 0.008188580408590829
    $ws.Range("D57").Value2 = 0.006794052106116533
    $ws.Range("E57").Value2 = -0.001912759505481643
    $ws.Range("D58").Value2 = 0.00558548104567229
    $ws.Range("E58").Value2 = -0.00102145045965274
    $ws.Range("D59").Value2 = 0.0054528330024528
    $ws.Range("E59").Value2 = -0.003400470834423075
    $ws.Range("D60").Value2 = 0.004943318714578654
    $ws.Range("E60").Value2 = 0.005642472428827983
    $ws.Range("D61").Value2 = 0.004810155610259447
    $ws.Range("E61").Value2 = -0.005164446860560057
    $ws.Range("D62").Value2 = 0.004944982758131587
    $ws.Range("E62").Value2 = -0.01610447880778776
    $ws.Range("D63").Value2 = 0.004202264652339305
    $ws.Range("E63").Value2 = 0.0007542615779152495
    $ws.Range("D64").Value2 = 0.004155354472180441
    $ws.Range("E64").Value2 = 0.02212051868802445
    $ws.Range("D65").Value2 = 0.00384536693031984
    $ws.Range("E65").Value2 = 0.008860863831190224
    $ws.Range("D66").Value2 = 0.003786412244444512
    $ws.Range("E66").Value2 = 0.01180311401305878
    $ws.Range("D67").Value2 = 0.003800913195405782
    $ws.Range("E67").Value2 = 0.01119519669772973
    $ws.Range("D68").Value2 = 0.003660063794675409
    $ws.Range("E68").Value2 = 0.003637190270516122
    $ws.Range("D69").Value2 = 0.003523017922066026
    $ws.Range("E69").Value2 = -0.01090868196131356
    $ws.Range("D70").Value2 = 0.00296532561132602
    $ws.Range("E70").Value2 = -0.002351557907113611
    $ws.Range("D71").Value2 = 0.00290767838824228
    $ws.Range("E71").Value2 = 0.02824673997465554
    $ws.Range("D72").Value2 = 0.002229144819491687
    $ws.Range("E72").Value2 = 0.03887101647619229
    $ws.Range("D73").Value2 = 0.00193227152563158
    $ws.Range("E73").Value2 = 0.02864465860159937
    $ws.Range("D74").Value2 = 0.001903626204470382
    $ws.Range("E74").Value2 = 0.02331050846046567
    $ws.Range("D75").Value2 = 0.001535357518172544
    $ws.Range("E75").Value2 = -0.005419075144508567
    $ws.Range("D76").Value2 = 0.00168147639015387
    $ws.Range("E76").Value2 = -0.006361922714420376
    $ws.Range("E77").Value2 = 0.009945449482728241

    $ws.Protect()
    Write-Output "Edit applied successfully"
} catch {
    Write-Output "ERROR: $_"
}